# Apply the "cfast source: removed unused variable NEUTRAL from CFAST data
# common block" edit to the VandV Comparison workbook.
#
# The underlying data change removes a no-longer-used NEUTRAL variable from
# the CFAST common block, which shifts/recomputes several sigma-M (K column)
# values in the comparison table (now derived as half of an updated total),
# and nudges a handful of offset (J column) values to match the refreshed
# CFAST output. Two previously-zero sigma-M cells become blank because the
# corresponding model no longer reports a value for those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: HGL Temperature Rise ---
$ws.Range("J3").Value2 = 1.1100000000000001
$ws.Range("K3").Formula = "=0.44/2"

# --- Row 4: HGL Depth ---
$ws.Range("J4").Value2 = 1.01
$ws.Range("K4").Formula = "=0.32/2"

# --- Row 5: Ceiling Jet Temp. Rise ---
$ws.Range("J5").Value2 = 1.25
$ws.Range("K5").Formula = "=0.53/2"

# --- Row 6: Plume Temperature Rise (only sigma M changes) ---
$ws.Range("K6").Formula = "=0.42/2"

# --- Row 8: Oxygen Concentration ---
$ws.Range("J8").Value2 = 1.03
$ws.Range("K8").Formula = "=0.63/2"

# --- Row 9: Smoke Concentration (only sigma M changes) ---
$ws.Range("K9").Formula = "=0.56/2"

# --- Row 10: Room Pressure Rise - sigma M no longer reported (was 0) ---
$ws.Range("K10").ClearContents()

# --- Row 11: Target Temperature Rise - sigma M no longer reported (was 0) ---
$ws.Range("K11").ClearContents()

# --- Row 13: Radiant Heat Flux (only sigma M changes) ---
$ws.Range("K13").Formula = "=1.29/2"

# --- Row 15: Total Heat Flux ---
$ws.Range("J15").Value2 = 0.99
$ws.Range("K15").Formula = "=0.99/2"

# --- Selection moved from F12 to the merged header J1:L1 ---
$null = $ws.Range("J1:L1").Select()

# --- Window height cosmetic change (persisted bookViews/workbookView) ---
# Best-effort: adjust the window's reported height to match the recorded
# view size. (Some headless hosts may not persist this purely cosmetic
# window-chrome value back into bookViews/workbookView.)
$win = $excel.Windows.Item(1)
$win.Height = 16440
